$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NroSiniestro in F3 with the new claim number (keep as text, preserving
# the existing General/quote-prefixed text style already applied to the cell)
$ws.Range("F3").Value = "'0420172010458  "

# Move the active cell selection to I5, matching the saved view state
$ws.Range("I5").Select()
